$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96, shifting existing rows 96-100 down to 97-101.
$ws.Rows.Item(96).Insert()

# Fill the static (unchanged) columns for the new row 96 by copying from the
# row directly below it (now row 97, formerly row 96) which shares the same
# market/region/product metadata.
$ws.Cells.Item(96, 1).Value = $ws.Cells.Item(97, 1).Value2
$ws.Cells.Item(96, 2).Value = $ws.Cells.Item(97, 2).Value2
$ws.Cells.Item(96, 3).Value = $ws.Cells.Item(97, 3).Value2
$ws.Cells.Item(96, 5).Value = $ws.Cells.Item(97, 5).Value2
$ws.Cells.Item(96, 6).Value = $ws.Cells.Item(97, 6).Value2
$ws.Cells.Item(96, 7).Value = $ws.Cells.Item(97, 7).Value2
$ws.Cells.Item(96, 8).Value = $ws.Cells.Item(97, 8).Value2
$ws.Cells.Item(96, 9).Value = $ws.Cells.Item(97, 9).Value2
$ws.Cells.Item(96, 10).Value = $ws.Cells.Item(97, 10).Value2
$ws.Cells.Item(96, 17).Value = $ws.Cells.Item(97, 17).Value2
$ws.Cells.Item(96, 18).Value = $ws.Cells.Item(97, 18).Value2
$ws.Cells.Item(96, 20).Value = $ws.Cells.Item(97, 20).Value2

# Set the new row's specific data values.
$ws.Cells.Item(96, 4).Value = 44706
$ws.Cells.Item(96, 11).Value = "Clemenuless"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 250
$ws.Cells.Item(96, 14).Value = 22000
$ws.Cells.Item(96, 15).Value = 23000
$ws.Cells.Item(96, 16).Value = 22500
$ws.Cells.Item(96, 19).Value = 1125
